$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier BSO"
$ws.Cells.Item($row, 3).Value = "VGO"

$dcell = $ws.Cells.Item($row, 4)
$dcell.NumberFormat = "@"
$dcell.Value = "2024-09-10"
$dcell.Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
